# "Fixed bug geek arrow rotation"
# Both existing issue rows (A1, A2) get a "Fixed" status marker added in
# column B, next to the issue description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Fixed"
$ws.Range("B2").Value = "Fixed"

# Leave the selection on the next empty cell below the new column, as in
# the saved workbook (activeCell/sqref = B3).
$ws.Range("B3").Select()
